$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 185 (shifts existing rows 185:295 down to 186:296,
# and the used range grows from A1:R295 to A1:R296).
$ws.Rows.Item(185).Insert()

# Populate the freshly inserted row 185 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R reuse the same values the (now shifted-down)
# row 186 carries; D,J,K,L,M,P get the new values for this record.
$ws.Cells.Item(185, 1).Value = 10
$ws.Cells.Item(185, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(185, 3).Value = "La Araucanía"
$ws.Cells.Item(185, 4).Value = "2022-01-21"
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 100112037
$ws.Cells.Item(185, 7).Value = "Cebollín"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 50
$ws.Cells.Item(185, 11).Value = 7000
$ws.Cells.Item(185, 12).Value = 7000
$ws.Cells.Item(185, 13).Value = 7000
$ws.Cells.Item(185, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(185, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(185, 16).Value = 583
$ws.Cells.Item(185, 17).Value = 12
$ws.Cells.Item(185, 18).Value = "Hortaliza"
